# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the file
# c9d153a9-612f-4e3d-8b76-142677743f2c.md across the Overview, zh-cn and
# de-de sheets to reflect a newly generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 5 is c9d153a9-612f-4e3d-8b76-142677743f2c.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-10-24 07:51:23"

# --- zh-cn sheet: row 5 is c9d153a9-612f-4e3d-8b76-142677743f2c.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-10-24 07:51:11"
$wsZhCn.Range("K5").Value = "2016-10-24 07:52:03"

# --- de-de sheet: row 5 is c9d153a9-612f-4e3d-8b76-142677743f2c.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-10-24 07:51:23"
$wsDeDe.Range("K5").Value = "2016-10-24 07:52:20"
